# Update the loading_percent values for the "380 kV" case
# Applies new computed values to columns B,D,E,F,G,H,J,K,M,N,O across rows 2-25
# Columns A, C, I, L are left unchanged (A = index, C/I/L = 0 constants)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 7.804822523898126
$ws.Cells.Item(3, 2).Value = 7.732920495321255
$ws.Cells.Item(4, 2).Value = 7.690262217059379
$ws.Cells.Item(5, 2).Value = 7.673272793920585
$ws.Cells.Item(6, 2).Value = 7.670476076672237
$ws.Cells.Item(7, 2).Value = 7.690031470496596
$ws.Cells.Item(8, 2).Value = 7.779733419642699
$ws.Cells.Item(9, 2).Value = 7.966569429422567
$ws.Cells.Item(10, 2).Value = 8.109272333978696
$ws.Cells.Item(11, 2).Value = 8.175103554152507
$ws.Cells.Item(12, 2).Value = 8.200141066719919
$ws.Cells.Item(13, 2).Value = 8.194744321704993
$ws.Cells.Item(14, 2).Value = 8.177161358651809
$ws.Cells.Item(15, 2).Value = 8.166404757907184
$ws.Cells.Item(16, 2).Value = 8.104986560330554
$ws.Cells.Item(17, 2).Value = 8.067526512580898
$ws.Cells.Item(18, 2).Value = 8.046067975407363
$ws.Cells.Item(19, 2).Value = 8.038818209064591
$ws.Cells.Item(20, 2).Value = 8.071505309732141
$ws.Cells.Item(21, 2).Value = 8.18232313435411
$ws.Cells.Item(22, 2).Value = 8.255371788018531
$ws.Cells.Item(23, 2).Value = 8.216334902907162
$ws.Cells.Item(24, 2).Value = 8.069706251641644
$ws.Cells.Item(25, 2).Value = 7.914980527065076
$ws.Cells.Item(2, 4).Value = 9.99258671121394
$ws.Cells.Item(3, 4).Value = 9.973345109851619
$ws.Cells.Item(4, 4).Value = 9.963112535756586
$ws.Cells.Item(5, 4).Value = 9.959343780021321
$ws.Cells.Item(6, 4).Value = 9.958742298300548
$ws.Cells.Item(7, 4).Value = 9.963060080759529
$ws.Cells.Item(8, 4).Value = 9.985625977456481
$ws.Cells.Item(9, 4).Value = 10.0422709407576
$ws.Cells.Item(10, 4).Value = 10.09121511838569
$ws.Cells.Item(11, 4).Value = 10.11501782948865
$ws.Cells.Item(12, 4).Value = 10.12424739060753
$ws.Cells.Item(13, 4).Value = 10.12225011769671
$ws.Cells.Item(14, 4).Value = 10.11577285257386
$ws.Cells.Item(15, 4).Value = 10.1118333122969
$ws.Cells.Item(16, 4).Value = 10.08969006558589
$ws.Cells.Item(17, 4).Value = 10.07649618446027
$ws.Cells.Item(18, 4).Value = 10.06905247814622
$ws.Cells.Item(19, 4).Value = 10.06655723449395
$ws.Cells.Item(20, 4).Value = 10.07788571834071
$ws.Cells.Item(21, 4).Value = 10.1176695626063
$ws.Cells.Item(22, 4).Value = 10.14492700762964
$ws.Cells.Item(23, 4).Value = 10.13026599681229
$ws.Cells.Item(24, 4).Value = 10.07725706852881
$ws.Cells.Item(25, 4).Value = 10.02564262174532
$ws.Cells.Item(2, 5).Value = 14.27630618581365
$ws.Cells.Item(3, 5).Value = 14.28493190878059
$ws.Cells.Item(4, 5).Value = 14.29265434284498
$ws.Cells.Item(5, 5).Value = 14.29641213072264
$ws.Cells.Item(6, 5).Value = 14.29707301939584
$ws.Cells.Item(7, 5).Value = 14.29270254759652
$ws.Cells.Item(8, 5).Value = 14.27877727602221
$ws.Cells.Item(9, 5).Value = 14.27068159016482
$ws.Cells.Item(10, 5).Value = 14.27638353981575
$ws.Cells.Item(11, 5).Value = 14.28148979853093
$ws.Cells.Item(12, 5).Value = 14.28378283088771
$ws.Cells.Item(13, 5).Value = 14.28327303303577
$ws.Cells.Item(14, 5).Value = 14.28167125624739
$ws.Cells.Item(15, 5).Value = 14.28073686711301
$ws.Cells.Item(16, 5).Value = 14.2761002343348
$ws.Cells.Item(17, 5).Value = 14.27389803189229
$ws.Cells.Item(18, 5).Value = 14.27286798976434
$ws.Cells.Item(19, 5).Value = 14.27255992626661
$ws.Cells.Item(20, 5).Value = 14.27410798514484
$ws.Cells.Item(21, 5).Value = 14.28213199887811
$ws.Cells.Item(22, 5).Value = 14.28947000608339
$ws.Cells.Item(23, 5).Value = 14.28536266309635
$ws.Cells.Item(24, 5).Value = 14.27401232992959
$ws.Cells.Item(25, 5).Value = 14.27082053064935
$ws.Cells.Item(2, 6).Value = 29.65268827905372
$ws.Cells.Item(3, 6).Value = 29.69423864898153
$ws.Cells.Item(4, 6).Value = 29.72641538546965
$ws.Cells.Item(5, 6).Value = 29.74120195313427
$ws.Cells.Item(6, 6).Value = 29.74375831684178
$ws.Cells.Item(7, 6).Value = 29.72660802568778
$ws.Cells.Item(8, 6).Value = 29.66563080244586
$ws.Cells.Item(9, 6).Value = 29.59898576335785
$ws.Cells.Item(10, 6).Value = 29.58232708838353
$ws.Cells.Item(11, 6).Value = 29.58175889120626
$ws.Cells.Item(12, 6).Value = 29.58255018996947
$ws.Cells.Item(13, 6).Value = 29.58233503886096
$ws.Cells.Item(14, 6).Value = 29.58180383212279
$ws.Cells.Item(15, 6).Value = 29.58160946616053
$ws.Cells.Item(16, 6).Value = 29.58250516131895
$ws.Cells.Item(17, 6).Value = 29.58484915640355
$ws.Cells.Item(18, 6).Value = 29.58685730247026
$ws.Cells.Item(19, 6).Value = 29.58765060829009
$ws.Cells.Item(20, 6).Value = 29.58453134171784
$ws.Cells.Item(21, 6).Value = 29.58193256104465
$ws.Cells.Item(22, 6).Value = 29.58609934369945
$ws.Cells.Item(23, 6).Value = 29.58333947214176
$ws.Cells.Item(24, 6).Value = 29.58467296806312
$ws.Cells.Item(25, 6).Value = 29.61134208029572
$ws.Cells.Item(2, 7).Value = 28.87574153429919
$ws.Cells.Item(3, 7).Value = 28.90457072843257
$ws.Cells.Item(4, 7).Value = 28.9308590992901
$ws.Cells.Item(5, 7).Value = 28.94372705021889
$ws.Cells.Item(6, 7).Value = 28.9459937906203
$ws.Cells.Item(7, 7).Value = 28.9310239208618
$ws.Cells.Item(8, 7).Value = 28.88389749846198
$ws.Cells.Item(9, 7).Value = 28.85974569627266
$ws.Cells.Item(10, 7).Value = 28.88371059087126
$ws.Cells.Item(11, 7).Value = 28.90366128779007
$ws.Cells.Item(12, 7).Value = 28.9125141896342
$ws.Cells.Item(13, 7).Value = 28.91054989651711
$ws.Cells.Item(14, 7).Value = 28.90436363103152
$ws.Cells.Item(15, 7).Value = 28.90074328291399
$ws.Cells.Item(16, 7).Value = 28.88258861183317
$ws.Cells.Item(17, 7).Value = 28.87376713311445
$ws.Cells.Item(18, 7).Value = 28.86954526689318
$ws.Cells.Item(19, 7).Value = 28.86826222417188
$ws.Cells.Item(20, 7).Value = 28.87461803291201
$ws.Cells.Item(21, 7).Value = 28.90614549257815
$ws.Cells.Item(22, 7).Value = 28.93431375458061
$ws.Cells.Item(23, 7).Value = 28.91858922788005
$ws.Cells.Item(24, 7).Value = 28.87423069407938
$ws.Cells.Item(25, 7).Value = 28.85895621975866
$ws.Cells.Item(2, 8).Value = 14.35245469603212
$ws.Cells.Item(3, 8).Value = 14.39237198742065
$ws.Cells.Item(4, 8).Value = 14.41901331321613
$ws.Cells.Item(5, 8).Value = 14.43040602921685
$ws.Cells.Item(6, 8).Value = 14.43233016547893
$ws.Cells.Item(7, 8).Value = 14.41916478840256
$ws.Cells.Item(8, 8).Value = 14.36577573420824
$ws.Cells.Item(9, 8).Value = 14.27799299543005
$ws.Cells.Item(10, 8).Value = 14.22380524790154
$ws.Cells.Item(11, 8).Value = 14.20139066571403
$ws.Cells.Item(12, 8).Value = 14.19322425924629
$ws.Cells.Item(13, 8).Value = 14.19496874162784
$ws.Cells.Item(14, 8).Value = 14.20071236695386
$ws.Cells.Item(15, 8).Value = 14.20427237364257
$ws.Cells.Item(16, 8).Value = 14.22531508341907
$ws.Cells.Item(17, 8).Value = 14.23879673753243
$ws.Cells.Item(18, 8).Value = 14.24676146238055
$ws.Cells.Item(19, 8).Value = 14.24949432280366
$ws.Cells.Item(20, 8).Value = 14.23733981312431
$ws.Cells.Item(21, 8).Value = 14.19901659921463
$ws.Cells.Item(22, 8).Value = 14.17584416478641
$ws.Cells.Item(23, 8).Value = 14.18804025529854
$ws.Cells.Item(24, 8).Value = 14.23799782223113
$ws.Cells.Item(25, 8).Value = 14.29993031100161
$ws.Cells.Item(2, 10).Value = 10.15113047360597
$ws.Cells.Item(3, 10).Value = 10.17403907347123
$ws.Cells.Item(4, 10).Value = 10.18925762855908
$ws.Cells.Item(5, 10).Value = 10.19574946070631
$ws.Cells.Item(6, 10).Value = 10.1968449580792
$ws.Cells.Item(7, 10).Value = 10.18934400445455
$ws.Cells.Item(8, 10).Value = 10.1587903344355
$ws.Cells.Item(9, 10).Value = 10.10800558502853
$ws.Cells.Item(10, 10).Value = 10.07624090870141
$ws.Cells.Item(11, 10).Value = 10.06299061220963
$ws.Cells.Item(12, 10).Value = 10.05814524636438
$ws.Cells.Item(13, 10).Value = 10.05918112676167
$ws.Cells.Item(14, 10).Value = 10.0625885307731
$ws.Cells.Item(15, 10).Value = 10.06469808601766
$ws.Cells.Item(16, 10).Value = 10.07713096309
$ws.Cells.Item(17, 10).Value = 10.08506518969733
$ws.Cells.Item(18, 10).Value = 10.0897416765977
$ws.Cells.Item(19, 10).Value = 10.09134445863859
$ws.Cells.Item(20, 10).Value = 10.08420889188226
$ws.Cells.Item(21, 10).Value = 10.06158302139885
$ws.Cells.Item(22, 10).Value = 10.04779949107544
$ws.Cells.Item(23, 10).Value = 10.05506426392952
$ws.Cells.Item(24, 10).Value = 10.0845956659091
$ws.Cells.Item(25, 10).Value = 10.12076869481245
$ws.Cells.Item(2, 11).Value = 9.212049532784579
$ws.Cells.Item(3, 11).Value = 8.890999866336644
$ws.Cells.Item(4, 11).Value = 8.688439382056902
$ws.Cells.Item(5, 11).Value = 8.604646497660166
$ws.Cells.Item(6, 11).Value = 8.590661081238171
$ws.Cells.Item(7, 11).Value = 8.687314209073076
$ws.Cells.Item(8, 11).Value = 9.102547304513541
$ws.Cells.Item(9, 11).Value = 9.869051635706001
$ws.Cells.Item(10, 11).Value = 10.44531280295914
$ws.Cells.Item(11, 11).Value = 10.71500938194179
$ws.Cells.Item(12, 11).Value = 10.8151428695635
$ws.Cells.Item(13, 11).Value = 10.79366656602625
$ws.Cells.Item(14, 11).Value = 10.72328758986536
$ws.Cells.Item(15, 11).Value = 10.67991767773957
$ws.Cells.Item(16, 11).Value = 10.4274102139228
$ws.Cells.Item(17, 11).Value = 10.26898940787663
$ws.Cells.Item(18, 11).Value = 10.17659012741652
$ws.Cells.Item(19, 11).Value = 10.14508674110429
$ws.Cells.Item(20, 11).Value = 10.28598631793947
$ws.Cells.Item(21, 11).Value = 10.74401399942017
$ws.Cells.Item(22, 11).Value = 11.03172195906216
$ws.Cells.Item(23, 11).Value = 10.87924235279221
$ws.Cells.Item(24, 11).Value = 10.27830612684158
$ws.Cells.Item(25, 11).Value = 9.667409365451956
$ws.Cells.Item(2, 13).Value = 14.57408802428979
$ws.Cells.Item(3, 13).Value = 14.45897909058559
$ws.Cells.Item(4, 13).Value = 14.38979437091715
$ws.Cells.Item(5, 13).Value = 14.36200042235972
$ws.Cells.Item(6, 13).Value = 14.35741007015228
$ws.Cells.Item(7, 13).Value = 14.38941788362753
$ws.Cells.Item(8, 13).Value = 14.53410328922454
$ws.Cells.Item(9, 13).Value = 14.8285906005432
$ws.Cells.Item(10, 13).Value = 15.05002750691108
$ws.Cells.Item(11, 13).Value = 15.15154566944166
$ws.Cells.Item(12, 13).Value = 15.19007364243109
$ws.Cells.Item(13, 13).Value = 15.18177260297482
$ws.Cells.Item(14, 13).Value = 15.15471383889391
$ws.Cells.Item(15, 13).Value = 15.13814984537367
$ws.Cells.Item(16, 13).Value = 15.0434064926616
$ws.Cells.Item(17, 13).Value = 14.98546596684284
$ws.Cells.Item(18, 13).Value = 14.95221569362914
$ws.Cells.Item(19, 13).Value = 14.9409715483544
$ws.Cells.Item(20, 13).Value = 14.99162621269159
$ws.Cells.Item(21, 13).Value = 15.16265956063375
$ws.Cells.Item(22, 13).Value = 15.27492384878444
$ws.Cells.Item(23, 13).Value = 15.21497111857753
$ws.Cells.Item(24, 13).Value = 14.98884097495642
$ws.Cells.Item(25, 13).Value = 14.74792367020266
$ws.Cells.Item(2, 14).Value = 18.86448358266832
$ws.Cells.Item(3, 14).Value = 18.91688675823854
$ws.Cells.Item(4, 14).Value = 18.95080494045076
$ws.Cells.Item(5, 14).Value = 18.96506605717068
$ws.Cells.Item(6, 14).Value = 18.96746066088844
$ws.Cells.Item(7, 14).Value = 18.95099549124346
$ws.Cells.Item(8, 14).Value = 18.88219118211504
$ws.Cells.Item(9, 14).Value = 18.76104410973932
$ws.Cells.Item(10, 14).Value = 18.68037039574579
$ws.Cells.Item(11, 14).Value = 18.64546530851532
$ws.Cells.Item(12, 14).Value = 18.63250458364431
$ws.Cells.Item(13, 14).Value = 18.6352844878946
$ws.Cells.Item(14, 14).Value = 18.64439387466098
$ws.Cells.Item(15, 14).Value = 18.65000709015912
$ws.Cells.Item(16, 14).Value = 18.68268754612594
$ws.Cells.Item(17, 14).Value = 18.70319475532114
$ws.Cells.Item(18, 14).Value = 18.71515885067271
$ws.Cells.Item(19, 14).Value = 18.71923872178153
$ws.Cells.Item(20, 14).Value = 18.70099425466379
$ws.Cells.Item(21, 14).Value = 18.64171125689397
$ws.Cells.Item(22, 14).Value = 18.6044644229994
$ws.Cells.Item(23, 14).Value = 18.62420697266625
$ws.Cells.Item(24, 14).Value = 18.70198855823707
$ws.Cells.Item(25, 14).Value = 18.79234949802778
$ws.Cells.Item(2, 15).Value = 21.85696387625436
$ws.Cells.Item(3, 15).Value = 21.91506029678743
$ws.Cells.Item(4, 15).Value = 21.95509983281876
$ws.Cells.Item(5, 15).Value = 21.97251313753484
$ws.Cells.Item(6, 15).Value = 21.97547080469691
$ws.Cells.Item(7, 15).Value = 21.95533023566775
$ws.Cells.Item(8, 15).Value = 21.87608804194454
$ws.Cells.Item(9, 15).Value = 21.75541736010543
$ws.Cells.Item(10, 15).Value = 21.68801010489561
$ws.Cells.Item(11, 15).Value = 21.66197353970986
$ws.Cells.Item(12, 15).Value = 21.6527805110616
$ws.Cells.Item(13, 15).Value = 21.65473073869618
$ws.Cells.Item(14, 15).Value = 21.66120386159218
$ws.Cells.Item(15, 15).Value = 21.66525565639401
$ws.Cells.Item(16, 15).Value = 21.68980486296823
$ws.Cells.Item(17, 15).Value = 21.70605097707901
$ws.Cells.Item(18, 15).Value = 21.71583080659941
$ws.Cells.Item(19, 15).Value = 21.71921685069237
$ws.Cells.Item(20, 15).Value = 21.70427647032023
$ws.Cells.Item(21, 15).Value = 21.65928445362287
$ws.Cells.Item(22, 15).Value = 21.63376453510871
$ws.Cells.Item(23, 15).Value = 21.64702921947961
$ws.Cells.Item(24, 15).Value = 21.70507735520086
$ws.Cells.Item(25, 15).Value = 21.78433565770271
